# Normalize "Raffaele de Rosa" name by syncing Hoja1 with the full
# file_path/file_name list already present on Hoja2 (rows 2-23, years
# 2002-2023), and refresh both sheets' selections.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Copy the full data block (A2:B23) from Hoja2 into Hoja1 so both sheets
# list every year 2002-2023 with matching file_path / file_name pairs.
# (Value2 is used because Value's getter is not reliable in this host.)
$ws1.Range("A2:B23").Value2 = $ws2.Range("A2:B23").Value2

# Update the stored selections to match the newly populated range.
$ws1.Activate()
$ws1.Range("A2:B23").Select()

$ws2.Activate()
$ws2.Range("A2:B23").Select()

$ws1.Activate()
